$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C values to match Bestand1 and highlight with a red fill
$updates = @{
    "C4"  = "Bestand 1 column 4"
    "C6"  = "Bestand 1 column 6"
    "C7"  = "Bestand 1 column 2"
    "C9"  = "Bestand 1 column 8"
    "C10" = "Bestand 1 column 11"
    "C11" = "Bestand 1 column 10"
    "C12" = "Bestand 1 column 9"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Interior.Color = 255
}
